$wb = $excel.ActiveWorkbook

# Locate the source sheet that the new sheet is cloned from.
$src = $wb.Worksheets.Item("optimized get_part_supserset")

# Match the pageSetup orientation that appears on both the source sheet and
# its clone after this edit.
$src.PageSetup.Orientation = 1

# Duplicate the sheet, placing the copy immediately after the source sheet.
$src.Copy([System.Reflection.Missing]::Value, $src)
$newSheet = $wb.Worksheets.Item("optimized get_part_supserset (2)")
$newSheet.Name = "Installed player hit_info hash"

# Update the hitVals block (B2:D2) and the %-improvement formula so it
# references the new sheet instead of the sheet it was copied from.
$newSheet.Range("B2").Value = 1.19
$newSheet.Range("C2").Value = 1.222
$newSheet.Range("D2").Value = 1.2
$newSheet.Range("F2").Formula = "=('optimized get_part_supserset'!E2-'Installed player hit_info hash'!E2)/'optimized get_part_supserset'!E2"

# Update the otherInfo block (B3:D3).
$newSheet.Range("B3").Value = 1.105
$newSheet.Range("C3").Value = 1.132
$newSheet.Range("D3").Value = 1.113
$newSheet.Range("H3").Value = "Installed hash table-based lookup of player hit Infos (hitVals and otherInfo)"

# Update the U block (B4:D4).
$newSheet.Range("B4").Value = 0.081
$newSheet.Range("C4").Value = 0.085
$newSheet.Range("D4").Value = 0.083
$newSheet.Range("H4").Value = "During setup, csv's of the player's hit info for the entire season"

# Notes column for the remaining rows (values/formulas are copied verbatim
# from the source sheet, so only the note text needs to change).
$newSheet.Range("H5").Value = "are generated, and then lazily loaded into a hash table"
$newSheet.Range("H6").Value = "during the simulation"

# Restore the old sheet's selection to the full used range, then make the
# new sheet the active / selected tab with H6 selected.
$src.Range("A1:H6").Select() | Out-Null
$newSheet.Select() | Out-Null
$newSheet.Range("H6").Select() | Out-Null
